$d = $word.ActiveDocument
$d.Content.Find.Execute("2024-11-23 Saturday", $true, $true, $false, $false, $false, $true, 1, $false, "2024-11-24 Sunday", 2) | Out-Null
$d.Content.Find.Execute("74-69=", $true, $true, $false, $false, $false, $true, 1, $false, "27+10=", 2) | Out-Null
$d.Content.Find.Execute("2+0=", $true, $true, $false, $false, $false, $true, 1, $false, "1+9=", 2) | Out-Null
$d.Content.Find.Execute("2+86=", $true, $true, $false, $false, $false, $true, 1, $false, "21+78=", 2) | Out-Null
$d.Content.Find.Execute("14+36=", $true, $true, $false, $false, $false, $true, 1, $false, "94+1=", 2) | Out-Null
$d.Content.Find.Execute("88-49=", $true, $true, $false, $false, $false, $true, 1, $false, "94-37=", 2) | Out-Null
$d.Content.Find.Execute("1+23=", $true, $true, $false, $false, $false, $true, 1, $false, "51-46=", 2) | Out-Null
$d.Content.Find.Execute("17+79=", $true, $true, $false, $false, $false, $true, 1, $false, "23-2=", 2) | Out-Null
$d.Content.Find.Execute("40-8=", $true, $true, $false, $false, $false, $true, 1, $false, "49+18=", 2) | Out-Null
$d.Content.Find.Execute("65-20=", $true, $true, $false, $false, $false, $true, 1, $false, "74-36=", 2) | Out-Null
$d.Content.Find.Execute("22+3=", $true, $true, $false, $false, $false, $true, 1, $false, "72-38=", 2) | Out-Null
$d.Content.Find.Execute("95-47=", $true, $true, $false, $false, $false, $true, 1, $false, "25-6=", 2) | Out-Null
$d.Content.Find.Execute("20+62=", $true, $true, $false, $false, $false, $true, 1, $false, "48-21=", 2) | Out-Null
$d.Content.Find.Execute("98-28=", $true, $true, $false, $false, $false, $true, 1, $false, "74-31=", 2) | Out-Null
$d.Content.Find.Execute("61+13=", $true, $true, $false, $false, $false, $true, 1, $false, "96-39=", 2) | Out-Null
$d.Content.Find.Execute("7+59=", $true, $true, $false, $false, $false, $true, 1, $false, "38+45=", 2) | Out-Null
$d.Content.Find.Execute("26-4=", $true, $true, $false, $false, $false, $true, 1, $false, "37+19=", 2) | Out-Null
$d.Content.Find.Execute("25+4=", $true, $true, $false, $false, $false, $true, 1, $false, "89-43=", 2) | Out-Null
$d.Content.Find.Execute("93-62=", $true, $true, $false, $false, $false, $true, 1, $false, "73-14=", 2) | Out-Null
$d.Content.Find.Execute("7+87=", $true, $true, $false, $false, $false, $true, 1, $false, "62-60=", 2) | Out-Null
$d.Content.Find.Execute("80-15=", $true, $true, $false, $false, $false, $true, 1, $false, "88-19=", 2) | Out-Null
$d.Content.Find.Execute("16+72=", $true, $true, $false, $false, $false, $true, 1, $false, "52+44=", 2) | Out-Null
$d.Content.Find.Execute("13+48=", $true, $true, $false, $false, $false, $true, 1, $false, "47+30=", 2) | Out-Null
$d.Content.Find.Execute("98-77=", $true, $true, $false, $false, $false, $true, 1, $false, "67-60=", 2) | Out-Null
$d.Content.Find.Execute("57-8=", $true, $true, $false, $false, $false, $true, 1, $false, "22+28=", 2) | Out-Null
$d.Content.Find.Execute("54-52=", $true, $true, $false, $false, $false, $true, 1, $false, "64+12=", 2) | Out-Null
$d.Content.Find.Execute("43+31=", $true, $true, $false, $false, $false, $true, 1, $false, "67-33=", 2) | Out-Null
$d.Content.Find.Execute("60+1=", $true, $true, $false, $false, $false, $true, 1, $false, "12+60=", 2) | Out-Null
$d.Content.Find.Execute("1+36=", $true, $true, $false, $false, $false, $true, 1, $false, "15+31=", 2) | Out-Null
$d.Content.Find.Execute("9+38=", $true, $true, $false, $false, $false, $true, 1, $false, "8+43=", 2) | Out-Null
$d.Content.Find.Execute("57-5=", $true, $true, $false, $false, $false, $true, 1, $false, "63-57=", 2) | Out-Null
$d.Content.Find.Execute("24+5=", $true, $true, $false, $false, $false, $true, 1, $false, "5+12=", 2) | Out-Null
$d.Content.Find.Execute("54+21=", $true, $true, $false, $false, $false, $true, 1, $false, "27-14=", 2) | Out-Null
$d.Content.Find.Execute("30-12=", $true, $true, $false, $false, $false, $true, 1, $false, "65-27=", 2) | Out-Null
$d.Content.Find.Execute("97-41=", $true, $true, $false, $false, $false, $true, 1, $false, "96-21=", 2) | Out-Null
$d.Content.Find.Execute("88-76=", $true, $true, $false, $false, $false, $true, 1, $false, "2+97=", 2) | Out-Null
$d.Content.Find.Execute("52+38=", $true, $true, $false, $false, $false, $true, 1, $false, "84-50=", 2) | Out-Null
$d.Content.Find.Execute("4+14=", $true, $true, $false, $false, $false, $true, 1, $false, "83-68=", 2) | Out-Null
$d.Content.Find.Execute("81+9=", $true, $true, $false, $false, $false, $true, 1, $false, "43-8=", 2) | Out-Null
$d.Content.Find.Execute("62-40=", $true, $true, $false, $false, $false, $true, 1, $false, "88-11=", 2) | Out-Null
$d.Content.Find.Execute("70-19=", $true, $true, $false, $false, $false, $true, 1, $false, "46-38=", 2) | Out-Null
$d.Content.Find.Execute("14-12=", $true, $true, $false, $false, $false, $true, 1, $false, "47-0=", 2) | Out-Null
$d.Content.Find.Execute("36+1=", $true, $true, $false, $false, $false, $true, 1, $false, "69+19=", 2) | Out-Null
$d.Content.Find.Execute("10+7=", $true, $true, $false, $false, $false, $true, 1, $false, "78-4=", 2) | Out-Null
$d.Content.Find.Execute("97-16=", $true, $true, $false, $false, $false, $true, 1, $false, "96-63=", 2) | Out-Null
$d.Content.Find.Execute("71+27=", $true, $true, $false, $false, $false, $true, 1, $false, "4+84=", 2) | Out-Null
$d.Content.Find.Execute("99-49=", $true, $true, $false, $false, $false, $true, 1, $false, "92+3=", 2) | Out-Null
$d.Content.Find.Execute("21+53=", $true, $true, $false, $false, $false, $true, 1, $false, "22+33=", 2) | Out-Null
$d.Content.Find.Execute("12+38=", $true, $true, $false, $false, $false, $true, 1, $false, "76+8=", 2) | Out-Null
$d.Content.Find.Execute("93-42=", $true, $true, $false, $false, $false, $true, 1, $false, "30+9=", 2) | Out-Null
$d.Content.Find.Execute("34-24=", $true, $true, $false, $false, $false, $true, 1, $false, "84-58=", 2) | Out-Null
$d.Content.Find.Execute("56+28=", $true, $true, $false, $false, $false, $true, 1, $false, "23-19=", 2) | Out-Null
$d.Content.Find.Execute("97-73=", $true, $true, $false, $false, $false, $true, 1, $false, "26+15=", 2) | Out-Null
$d.Content.Find.Execute("39-9=", $true, $true, $false, $false, $false, $true, 1, $false, "25+19=", 2) | Out-Null
$d.Content.Find.Execute("79-43=", $true, $true, $false, $false, $false, $true, 1, $false, "73+2=", 2) | Out-Null
$d.Content.Find.Execute("63-9=", $true, $true, $false, $false, $false, $true, 1, $false, "43+30=", 2) | Out-Null
$d.Content.Find.Execute("76-59=", $true, $true, $false, $false, $false, $true, 1, $false, "17+18=", 2) | Out-Null
$d.Content.Find.Execute("17+8=", $true, $true, $false, $false, $false, $true, 1, $false, "88+2=", 2) | Out-Null
$d.Content.Find.Execute("33-16=", $true, $true, $false, $false, $false, $true, 1, $false, "90-74=", 2) | Out-Null
$d.Content.Find.Execute("7+46=", $true, $true, $false, $false, $false, $true, 1, $false, "34+26=", 2) | Out-Null
$d.Content.Find.Execute("48+49=", $true, $true, $false, $false, $false, $true, 1, $false, "89-83=", 2) | Out-Null
$d.Content.Find.Execute("77-1=", $true, $true, $false, $false, $false, $true, 1, $false, "57+5=", 2) | Out-Null
$d.Content.Find.Execute("24+49=", $true, $true, $false, $false, $false, $true, 1, $false, "7+26=", 2) | Out-Null
$d.Content.Find.Execute("73-9=", $true, $true, $false, $false, $false, $true, 1, $false, "83+14=", 2) | Out-Null
$d.Content.Find.Execute("60-10=", $true, $true, $false, $false, $false, $true, 1, $false, "11+80=", 2) | Out-Null
$d.Content.Find.Execute("58-47=", $true, $true, $false, $false, $false, $true, 1, $false, "17+63=", 2) | Out-Null
$d.Content.Find.Execute("33+4=", $true, $true, $false, $false, $false, $true, 1, $false, "62-32=", 2) | Out-Null
$d.Content.Find.Execute("82-39=", $true, $true, $false, $false, $false, $true, 1, $false, "17-0=", 2) | Out-Null
$d.Content.Find.Execute("27+71=", $true, $true, $false, $false, $false, $true, 1, $false, "10+20=", 2) | Out-Null
$d.Content.Find.Execute("69+30=", $true, $true, $false, $false, $false, $true, 1, $false, "98-15=", 2) | Out-Null
$d.Content.Find.Execute("8+29=", $true, $true, $false, $false, $false, $true, 1, $false, "88+1=", 2) | Out-Null
$d.Content.Find.Execute("77-28=", $true, $true, $false, $false, $false, $true, 1, $false, "91-34=", 2) | Out-Null
$d.Content.Find.Execute("96+1=", $true, $true, $false, $false, $false, $true, 1, $false, "42-38=", 2) | Out-Null
$d.Content.Find.Execute("38-37=", $true, $true, $false, $false, $false, $true, 1, $false, "21+54=", 2) | Out-Null
$d.Content.Find.Execute("55-32=", $true, $true, $false, $false, $false, $true, 1, $false, "67-2=", 2) | Out-Null
$d.Content.Find.Execute("96-36=", $true, $true, $false, $false, $false, $true, 1, $false, "43+32=", 2) | Out-Null
$d.Content.Find.Execute("5+62=", $true, $true, $false, $false, $false, $true, 1, $false, "71+15=", 2) | Out-Null
$d.Content.Find.Execute("1+17=", $true, $true, $false, $false, $false, $true, 1, $false, "99-92=", 2) | Out-Null
$d.Content.Find.Execute("17+73=", $true, $true, $false, $false, $false, $true, 1, $false, "74-4=", 2) | Out-Null
$d.Content.Find.Execute("15+76=", $true, $true, $false, $false, $false, $true, 1, $false, "5+41=", 2) | Out-Null
$d.Content.Find.Execute("94-90=", $true, $true, $false, $false, $false, $true, 1, $false, "42+32=", 2) | Out-Null
$d.Content.Find.Execute("19-8=", $true, $true, $false, $false, $false, $true, 1, $false, "83-15=", 2) | Out-Null
$d.Content.Find.Execute("31+30=", $true, $true, $false, $false, $false, $true, 1, $false, "9+50=", 2) | Out-Null
$d.Content.Find.Execute("9+67=", $true, $true, $false, $false, $false, $true, 1, $false, "29+58=", 2) | Out-Null
$d.Content.Find.Execute("13-12=", $true, $true, $false, $false, $false, $true, 1, $false, "88-23=", 2) | Out-Null
$d.Content.Find.Execute("4+32=", $true, $true, $false, $false, $false, $true, 1, $false, "45-27=", 2) | Out-Null
$d.Content.Find.Execute("85-27=", $true, $true, $false, $false, $false, $true, 1, $false, "49+18=", 2) | Out-Null
$d.Content.Find.Execute("26+62=", $true, $true, $false, $false, $false, $true, 1, $false, "43+56=", 2) | Out-Null
$d.Content.Find.Execute("43-5=", $true, $true, $false, $false, $false, $true, 1, $false, "69-1=", 2) | Out-Null
$d.Content.Find.Execute("61+22=", $true, $true, $false, $false, $false, $true, 1, $false, "59-48=", 2) | Out-Null
$d.Content.Find.Execute("97-76=", $true, $true, $false, $false, $false, $true, 1, $false, "29+39=", 2) | Out-Null
$d.Content.Find.Execute("89+10=", $true, $true, $false, $false, $false, $true, 1, $false, "51+18=", 2) | Out-Null
$d.Content.Find.Execute("45+51=", $true, $true, $false, $false, $false, $true, 1, $false, "4+29=", 2) | Out-Null
$d.Content.Find.Execute("55+14=", $true, $true, $false, $false, $false, $true, 1, $false, "55-47=", 2) | Out-Null
$d.Content.Find.Execute("47+44=", $true, $true, $false, $false, $false, $true, 1, $false, "76-31=", 2) | Out-Null
$d.Content.Find.Execute("86-27=", $true, $true, $false, $false, $false, $true, 1, $false, "42-8=", 2) | Out-Null
$d.Content.Find.Execute("97-75=", $true, $true, $false, $false, $false, $true, 1, $false, "76-0=", 2) | Out-Null
$d.Content.Find.Execute("65+19=", $true, $true, $false, $false, $false, $true, 1, $false, "38+43=", 2) | Out-Null
$d.Content.Find.Execute("2+41=", $true, $true, $false, $false, $false, $true, 1, $false, "77+16=", 2) | Out-Null
$d.Content.Find.Execute("37+52=", $true, $true, $false, $false, $false, $true, 1, $false, "57+38=", 2) | Out-Null
$d.Content.Find.Execute("18+1=", $true, $true, $false, $false, $false, $true, 1, $false, "87-37=", 2) | Out-Null
